$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed new shared strings in the exact target order via a scratch cell ---
$ws.Range("ZZ1").Value = 'Machine learning digit recognition'
$ws.Range("ZZ1").Value = 'find jobs'
$ws.Range("ZZ1").ClearContents()

# ---- Block starting row 229 ----
$ws.Range("A217").Copy()
$ws.Range("A229").PasteSpecial(-4122)
$ws.Range("A229").Value = ' Date'
$ws.Range("B217").Copy()
$ws.Range("B229").PasteSpecial(-4122)
$ws.Range("B229").Value = 45287
$ws.Range("C217").Copy()
$ws.Range("C229").PasteSpecial(-4122)
$ws.Range("C229").Value = 'Total Time '
$ws.Range("D217").Copy()
$ws.Range("D229").PasteSpecial(-4122)
$ws.Range("D229").Value = 9
$ws.Range("E217").Copy()
$ws.Range("E229").PasteSpecial(-4122)
$ws.Range("E229").Value = 'Pay'
$ws.Range("F217").Copy()
$ws.Range("F229").PasteSpecial(-4122)
$ws.Range("F229").Value = 'ZENBUSINESS $324.00'
$ws.Range("G217").Copy()
$ws.Range("G229").PasteSpecial(-4122)
$ws.Range("G229").Value = 'ZENBUSINESS $199.00'
$ws.Range("A218").Copy()
$ws.Range("A230").PasteSpecial(-4122)
$ws.Range("A230").Value = 'Time'
$ws.Range("B218").Copy()
$ws.Range("B230").PasteSpecial(-4122)
$ws.Range("B230").Value = 'Task Description'
$ws.Range("C218").Copy()
$ws.Range("C230").PasteSpecial(-4122)
$ws.Range("C230").Value = 'Type'
$ws.Range("E218").Copy()
$ws.Range("E230").PasteSpecial(-4122)
$ws.Range("E230").Value = 'Income'
$ws.Range("F218").Copy()
$ws.Range("F230").PasteSpecial(-4122)
$ws.Range("F230").Value = 0
$ws.Range("A219").Copy()
$ws.Range("A231").PasteSpecial(-4122)
$ws.Range("A231").Value = 0.4166666666666667
$ws.Range("B219").Copy()
$ws.Range("B231").PasteSpecial(-4122)
$ws.Range("B231").Value = 'Math simulation on Python'
$ws.Range("C219").Copy()
$ws.Range("C231").PasteSpecial(-4122)
$ws.Range("C231").Value = 'design'
$ws.Range("E219").Copy()
$ws.Range("E231").PasteSpecial(-4122)
$ws.Range("A220").Copy()
$ws.Range("A232").PasteSpecial(-4122)
$ws.Range("A232").Value = 0.4583333333333333
$ws.Range("B220").Copy()
$ws.Range("B232").PasteSpecial(-4122)
$ws.Range("B232").Value = 'Math simulation on Python'
$ws.Range("C220").Copy()
$ws.Range("C232").PasteSpecial(-4122)
$ws.Range("C232").Value = 'design'
$ws.Range("E220").Copy()
$ws.Range("E232").PasteSpecial(-4122)
$ws.Range("A221").Copy()
$ws.Range("A233").PasteSpecial(-4122)
$ws.Range("A233").Value = 0.5
$ws.Range("B221").Copy()
$ws.Range("B233").PasteSpecial(-4122)
$ws.Range("B233").Value = 'Math simulation on Python'
$ws.Range("C221").Copy()
$ws.Range("C233").PasteSpecial(-4122)
$ws.Range("C233").Value = 'design'
$ws.Range("E221").Copy()
$ws.Range("E233").PasteSpecial(-4122)
$ws.Range("A222").Copy()
$ws.Range("A234").PasteSpecial(-4122)
$ws.Range("A234").Value = '13:00 PM'
$ws.Range("B222").Copy()
$ws.Range("B234").PasteSpecial(-4122)
$ws.Range("B234").Value = 'Math simulation on Python'
$ws.Range("C222").Copy()
$ws.Range("C234").PasteSpecial(-4122)
$ws.Range("C234").Value = 'design'
$ws.Range("E222").Copy()
$ws.Range("E234").PasteSpecial(-4122)
$ws.Range("A223").Copy()
$ws.Range("A235").PasteSpecial(-4122)
$ws.Range("A235").Value = '14:000 PM'
$ws.Range("B223").Copy()
$ws.Range("B235").PasteSpecial(-4122)
$ws.Range("B235").Value = 'find jobs'
$ws.Range("C223").Copy()
$ws.Range("C235").PasteSpecial(-4122)
$ws.Range("C235").Value = 'Support'
$ws.Range("E223").Copy()
$ws.Range("E235").PasteSpecial(-4122)
$ws.Range("A224").Copy()
$ws.Range("A236").PasteSpecial(-4122)
$ws.Range("A236").Value = '15:00 PM'
$ws.Range("B224").Copy()
$ws.Range("B236").PasteSpecial(-4122)
$ws.Range("B236").Value = 'find jobs'
$ws.Range("C224").Copy()
$ws.Range("C236").PasteSpecial(-4122)
$ws.Range("C236").Value = 'Support'
$ws.Range("E224").Copy()
$ws.Range("E236").PasteSpecial(-4122)
$ws.Range("A225").Copy()
$ws.Range("A237").PasteSpecial(-4122)
$ws.Range("A237").Value = '16:00 PM'
$ws.Range("B225").Copy()
$ws.Range("B237").PasteSpecial(-4122)
$ws.Range("B237").Value = 'find jobs'
$ws.Range("C225").Copy()
$ws.Range("C237").PasteSpecial(-4122)
$ws.Range("C237").Value = 'Support'
$ws.Range("E225").Copy()
$ws.Range("E237").PasteSpecial(-4122)
$ws.Range("A226").Copy()
$ws.Range("A238").PasteSpecial(-4122)
$ws.Range("A238").Value = '17:00 PM'
$ws.Range("B226").Copy()
$ws.Range("B238").PasteSpecial(-4122)
$ws.Range("B238").Value = 'Math simulation on Python'
$ws.Range("C226").Copy()
$ws.Range("C238").PasteSpecial(-4122)
$ws.Range("C238").Value = 'design'
$ws.Range("E226").Copy()
$ws.Range("E238").PasteSpecial(-4122)
$ws.Range("A227").Copy()
$ws.Range("A239").PasteSpecial(-4122)
$ws.Range("A239").Value = '18:00 PM'
$ws.Range("B227").Copy()
$ws.Range("B239").PasteSpecial(-4122)
$ws.Range("B239").Value = 'Math simulation on Python'
$ws.Range("C227").Copy()
$ws.Range("C239").PasteSpecial(-4122)
$ws.Range("C239").Value = 'design'
$ws.Range("E227").Copy()
$ws.Range("E239").PasteSpecial(-4122)
$ws.Range("A228").Copy()
$ws.Range("A240").PasteSpecial(-4122)
$ws.Range("A240").Value = '19:00 PM'
$ws.Range("B228").Copy()
$ws.Range("B240").PasteSpecial(-4122)
$ws.Range("B240").Value = 'Math simulation on Python'
$ws.Range("C228").Copy()
$ws.Range("C240").PasteSpecial(-4122)
$ws.Range("C240").Value = 'design'
$ws.Range("E228").Copy()
$ws.Range("E240").PasteSpecial(-4122)

# ---- Block starting row 241 ----
$ws.Range("A217").Copy()
$ws.Range("A241").PasteSpecial(-4122)
$ws.Range("A241").Value = ' Date'
$ws.Range("B217").Copy()
$ws.Range("B241").PasteSpecial(-4122)
$ws.Range("B241").Value = 45288
$ws.Range("C217").Copy()
$ws.Range("C241").PasteSpecial(-4122)
$ws.Range("C241").Value = 'Total Time '
$ws.Range("D217").Copy()
$ws.Range("D241").PasteSpecial(-4122)
$ws.Range("D241").Value = 9
$ws.Range("E217").Copy()
$ws.Range("E241").PasteSpecial(-4122)
$ws.Range("E241").Value = 'Pay'
$ws.Range("F217").Copy()
$ws.Range("F241").PasteSpecial(-4122)
$ws.Range("F241").Value = 'ZENBUSINESS $324.00'
$ws.Range("G217").Copy()
$ws.Range("G241").PasteSpecial(-4122)
$ws.Range("G241").Value = 'ZENBUSINESS $199.00'
$ws.Range("A218").Copy()
$ws.Range("A242").PasteSpecial(-4122)
$ws.Range("A242").Value = 'Time'
$ws.Range("B218").Copy()
$ws.Range("B242").PasteSpecial(-4122)
$ws.Range("B242").Value = 'Task Description'
$ws.Range("C218").Copy()
$ws.Range("C242").PasteSpecial(-4122)
$ws.Range("C242").Value = 'Type'
$ws.Range("E218").Copy()
$ws.Range("E242").PasteSpecial(-4122)
$ws.Range("E242").Value = 'Income'
$ws.Range("F218").Copy()
$ws.Range("F242").PasteSpecial(-4122)
$ws.Range("F242").Value = 0
$ws.Range("A219").Copy()
$ws.Range("A243").PasteSpecial(-4122)
$ws.Range("A243").Value = 0.4166666666666667
$ws.Range("B219").Copy()
$ws.Range("B243").PasteSpecial(-4122)
$ws.Range("B243").Value = 'Machine learning digit recognition'
$ws.Range("C219").Copy()
$ws.Range("C243").PasteSpecial(-4122)
$ws.Range("C243").Value = 'design'
$ws.Range("E219").Copy()
$ws.Range("E243").PasteSpecial(-4122)
$ws.Range("A220").Copy()
$ws.Range("A244").PasteSpecial(-4122)
$ws.Range("A244").Value = 0.4583333333333333
$ws.Range("B220").Copy()
$ws.Range("B244").PasteSpecial(-4122)
$ws.Range("B244").Value = 'Machine learning digit recognition'
$ws.Range("C220").Copy()
$ws.Range("C244").PasteSpecial(-4122)
$ws.Range("C244").Value = 'design'
$ws.Range("E220").Copy()
$ws.Range("E244").PasteSpecial(-4122)
$ws.Range("A221").Copy()
$ws.Range("A245").PasteSpecial(-4122)
$ws.Range("A245").Value = 0.5
$ws.Range("B221").Copy()
$ws.Range("B245").PasteSpecial(-4122)
$ws.Range("B245").Value = 'Machine learning digit recognition'
$ws.Range("C221").Copy()
$ws.Range("C245").PasteSpecial(-4122)
$ws.Range("C245").Value = 'design'
$ws.Range("E221").Copy()
$ws.Range("E245").PasteSpecial(-4122)
$ws.Range("A222").Copy()
$ws.Range("A246").PasteSpecial(-4122)
$ws.Range("A246").Value = '13:00 PM'
$ws.Range("B222").Copy()
$ws.Range("B246").PasteSpecial(-4122)
$ws.Range("B246").Value = 'Machine learning digit recognition'
$ws.Range("C222").Copy()
$ws.Range("C246").PasteSpecial(-4122)
$ws.Range("C246").Value = 'design'
$ws.Range("E222").Copy()
$ws.Range("E246").PasteSpecial(-4122)
$ws.Range("A223").Copy()
$ws.Range("A247").PasteSpecial(-4122)
$ws.Range("A247").Value = '14:000 PM'
$ws.Range("B223").Copy()
$ws.Range("B247").PasteSpecial(-4122)
$ws.Range("B247").Value = 'Machine learning digit recognition'
$ws.Range("C223").Copy()
$ws.Range("C247").PasteSpecial(-4122)
$ws.Range("C247").Value = 'design'
$ws.Range("E223").Copy()
$ws.Range("E247").PasteSpecial(-4122)
$ws.Range("A224").Copy()
$ws.Range("A248").PasteSpecial(-4122)
$ws.Range("A248").Value = '15:00 PM'
$ws.Range("B224").Copy()
$ws.Range("B248").PasteSpecial(-4122)
$ws.Range("B248").Value = 'find jobs'
$ws.Range("C224").Copy()
$ws.Range("C248").PasteSpecial(-4122)
$ws.Range("C248").Value = 'Support'
$ws.Range("E224").Copy()
$ws.Range("E248").PasteSpecial(-4122)
$ws.Range("A225").Copy()
$ws.Range("A249").PasteSpecial(-4122)
$ws.Range("A249").Value = '16:00 PM'
$ws.Range("B225").Copy()
$ws.Range("B249").PasteSpecial(-4122)
$ws.Range("B249").Value = 'find jobs'
$ws.Range("C225").Copy()
$ws.Range("C249").PasteSpecial(-4122)
$ws.Range("C249").Value = 'Support'
$ws.Range("E225").Copy()
$ws.Range("E249").PasteSpecial(-4122)
$ws.Range("A226").Copy()
$ws.Range("A250").PasteSpecial(-4122)
$ws.Range("A250").Value = '17:00 PM'
$ws.Range("B226").Copy()
$ws.Range("B250").PasteSpecial(-4122)
$ws.Range("B250").Value = 'find jobs'
$ws.Range("C226").Copy()
$ws.Range("C250").PasteSpecial(-4122)
$ws.Range("C250").Value = 'Support'
$ws.Range("E226").Copy()
$ws.Range("E250").PasteSpecial(-4122)
$ws.Range("A227").Copy()
$ws.Range("A251").PasteSpecial(-4122)
$ws.Range("A251").Value = '18:00 PM'
$ws.Range("B227").Copy()
$ws.Range("B251").PasteSpecial(-4122)
$ws.Range("B251").Value = 'Machine learning digit recognition'
$ws.Range("C227").Copy()
$ws.Range("C251").PasteSpecial(-4122)
$ws.Range("C251").Value = 'design'
$ws.Range("E227").Copy()
$ws.Range("E251").PasteSpecial(-4122)
$ws.Range("A228").Copy()
$ws.Range("A252").PasteSpecial(-4122)
$ws.Range("A252").Value = '19:00 PM'
$ws.Range("B228").Copy()
$ws.Range("B252").PasteSpecial(-4122)
$ws.Range("B252").Value = 'Machine learning digit recognition'
$ws.Range("C228").Copy()
$ws.Range("C252").PasteSpecial(-4122)
$ws.Range("C252").Value = 'design'
$ws.Range("E228").Copy()
$ws.Range("E252").PasteSpecial(-4122)

# ---- Block starting row 253 ----
$ws.Range("A217").Copy()
$ws.Range("A253").PasteSpecial(-4122)
$ws.Range("A253").Value = ' Date'
$ws.Range("B217").Copy()
$ws.Range("B253").PasteSpecial(-4122)
$ws.Range("B253").Value = 45289
$ws.Range("C217").Copy()
$ws.Range("C253").PasteSpecial(-4122)
$ws.Range("C253").Value = 'Total Time '
$ws.Range("D217").Copy()
$ws.Range("D253").PasteSpecial(-4122)
$ws.Range("D253").Value = 4
$ws.Range("E217").Copy()
$ws.Range("E253").PasteSpecial(-4122)
$ws.Range("E253").Value = 'Pay'
$ws.Range("F217").Copy()
$ws.Range("F253").PasteSpecial(-4122)
$ws.Range("F253").Value = 'ZENBUSINESS $324.00'
$ws.Range("G217").Copy()
$ws.Range("G253").PasteSpecial(-4122)
$ws.Range("G253").Value = 'ZENBUSINESS $199.00'
$ws.Range("A218").Copy()
$ws.Range("A254").PasteSpecial(-4122)
$ws.Range("A254").Value = 'Time'
$ws.Range("B218").Copy()
$ws.Range("B254").PasteSpecial(-4122)
$ws.Range("B254").Value = 'Task Description'
$ws.Range("C218").Copy()
$ws.Range("C254").PasteSpecial(-4122)
$ws.Range("C254").Value = 'Type'
$ws.Range("E218").Copy()
$ws.Range("E254").PasteSpecial(-4122)
$ws.Range("E254").Value = 'Income'
$ws.Range("F218").Copy()
$ws.Range("F254").PasteSpecial(-4122)
$ws.Range("F254").Value = 0
$ws.Range("A219").Copy()
$ws.Range("A255").PasteSpecial(-4122)
$ws.Range("A255").Value = 0.4166666666666667
$ws.Range("B219").Copy()
$ws.Range("B255").PasteSpecial(-4122)
$ws.Range("B255").Value = 'Machine learning digit recognition'
$ws.Range("C219").Copy()
$ws.Range("C255").PasteSpecial(-4122)
$ws.Range("C255").Value = 'design'
$ws.Range("E219").Copy()
$ws.Range("E255").PasteSpecial(-4122)
$ws.Range("A220").Copy()
$ws.Range("A256").PasteSpecial(-4122)
$ws.Range("A256").Value = 0.4583333333333333
$ws.Range("B220").Copy()
$ws.Range("B256").PasteSpecial(-4122)
$ws.Range("B256").Value = 'Machine learning digit recognition'
$ws.Range("C220").Copy()
$ws.Range("C256").PasteSpecial(-4122)
$ws.Range("C256").Value = 'design'
$ws.Range("E220").Copy()
$ws.Range("E256").PasteSpecial(-4122)
$ws.Range("A221").Copy()
$ws.Range("A257").PasteSpecial(-4122)
$ws.Range("A257").Value = 0.5
$ws.Range("B221").Copy()
$ws.Range("B257").PasteSpecial(-4122)
$ws.Range("B257").Value = 'Machine learning digit recognition'
$ws.Range("C221").Copy()
$ws.Range("C257").PasteSpecial(-4122)
$ws.Range("C257").Value = 'design'
$ws.Range("E221").Copy()
$ws.Range("E257").PasteSpecial(-4122)
$ws.Range("A222").Copy()
$ws.Range("A258").PasteSpecial(-4122)
$ws.Range("A258").Value = '13:00 PM'
$ws.Range("B222").Copy()
$ws.Range("B258").PasteSpecial(-4122)
$ws.Range("B258").Value = 'Machine learning digit recognition'
$ws.Range("C222").Copy()
$ws.Range("C258").PasteSpecial(-4122)
$ws.Range("C258").Value = 'design'
$ws.Range("E222").Copy()
$ws.Range("E258").PasteSpecial(-4122)
$ws.Range("A223").Copy()
$ws.Range("A259").PasteSpecial(-4122)
$ws.Range("A259").Value = '14:000 PM'
$ws.Range("B223").Copy()
$ws.Range("B259").PasteSpecial(-4122)
$ws.Range("B259").Value = 'Machine learning digit recognition'
$ws.Range("C223").Copy()
$ws.Range("C259").PasteSpecial(-4122)
$ws.Range("C259").Value = 'design'
$ws.Range("E223").Copy()
$ws.Range("E259").PasteSpecial(-4122)
$ws.Range("A224").Copy()
$ws.Range("A260").PasteSpecial(-4122)
$ws.Range("A260").Value = '15:00 PM'
$ws.Range("B224").Copy()
$ws.Range("B260").PasteSpecial(-4122)
$ws.Range("C224").Copy()
$ws.Range("C260").PasteSpecial(-4122)
$ws.Range("E224").Copy()
$ws.Range("E260").PasteSpecial(-4122)
$ws.Range("A225").Copy()
$ws.Range("A261").PasteSpecial(-4122)
$ws.Range("A261").Value = '16:00 PM'
$ws.Range("B225").Copy()
$ws.Range("B261").PasteSpecial(-4122)
$ws.Range("C225").Copy()
$ws.Range("C261").PasteSpecial(-4122)
$ws.Range("E225").Copy()
$ws.Range("E261").PasteSpecial(-4122)
$ws.Range("A226").Copy()
$ws.Range("A262").PasteSpecial(-4122)
$ws.Range("A262").Value = '17:00 PM'
$ws.Range("B226").Copy()
$ws.Range("B262").PasteSpecial(-4122)
$ws.Range("C226").Copy()
$ws.Range("C262").PasteSpecial(-4122)
$ws.Range("E226").Copy()
$ws.Range("E262").PasteSpecial(-4122)
$ws.Range("A227").Copy()
$ws.Range("A263").PasteSpecial(-4122)
$ws.Range("A263").Value = '18:00 PM'
$ws.Range("B227").Copy()
$ws.Range("B263").PasteSpecial(-4122)
$ws.Range("C227").Copy()
$ws.Range("C263").PasteSpecial(-4122)
$ws.Range("E227").Copy()
$ws.Range("E263").PasteSpecial(-4122)
$ws.Range("A228").Copy()
$ws.Range("A264").PasteSpecial(-4122)
$ws.Range("A264").Value = '19:00 PM'
$ws.Range("B228").Copy()
$ws.Range("B264").PasteSpecial(-4122)
$ws.Range("C228").Copy()
$ws.Range("C264").PasteSpecial(-4122)
$ws.Range("E228").Copy()
$ws.Range("E264").PasteSpecial(-4122)

$ws.Range("C262").Select()
Write-Host "done"